$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 141.25
$ws.Range("I33").Value = 144.375
$ws.Range("K33").Value = 144.375
$ws.Range("M33").Value = 84.625

$ws.Range("H100").Value = 13891598
$ws.Range("I100").Value = 23811496
$ws.Range("J100").Value = 3740
$ws.Range("K100").Value = 23811496
$ws.Range("L100").Value = 3740
$ws.Range("M100").Value = -23810955
$ws.Range("N100").Value = -4822

$ws.Range("H116").Value = 5500
$ws.Range("I116").Value = 5500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 5500
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -2058

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28275.32
$ws.Range("I32").Value = 6001.641
$ws.Range("K32").Value = 6001.641
$ws.Range("M32").Value = -5714.641

$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 30000
$ws.Range("N87").Value = -32496

$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 90000
$ws.Range("N90").Value = -102480

$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774

$ws.Range("H132").Value = 2914.525
$ws.Range("I132").Value = 2326.1936
$ws.Range("K132").Value = 6978.5808
$ws.Range("M132").Value = -4448.5808

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2805.44
$ws.Range("I134").Value = 2074.7896
$ws.Range("J134").Value = 5119.1665
$ws.Range("K134").Value = 6224.3688
$ws.Range("L134").Value = 15357.4995
$ws.Range("M134").Value = -3689.3688
$ws.Range("N134").Value = -20427.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4718.3223
$ws.Range("I31").Value = 1512.9706
$ws.Range("J31").Value = 9077.6
$ws.Range("K31").Value = 1512.9706
$ws.Range("L31").Value = 9077.6
$ws.Range("M31").Value = -1217.9706
$ws.Range("N31").Value = -9667.6

$ws.Range("H34").Value = 4718.3223
$ws.Range("I34").Value = 1512.9706
$ws.Range("J34").Value = 9077.6
$ws.Range("K34").Value = 1512.9706
$ws.Range("L34").Value = 9077.6
$ws.Range("M34").Value = -1310.9706
$ws.Range("N34").Value = -9481.6

$ws.Range("H58").Value = 3555.25
$ws.Range("I58").Value = 2200
$ws.Range("K58").Value = 2200
$ws.Range("M58").Value = -1997

$ws.Range("H132").Value = 2746.9285
$ws.Range("I132").Value = 1711.5555
$ws.Range("K132").Value = 5134.666499999999
$ws.Range("M132").Value = -2604.666499999999

$ws.Range("H136").Value = 3555.25
$ws.Range("I136").Value = 2200
$ws.Range("K136").Value = 6600
$ws.Range("M136").Value = -4050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5051846.5
$ws.Range("J131").Value = 5953896.5
$ws.Range("L131").Value = 17861689.5
$ws.Range("N131").Value = -17871769.5

$ws.Range("H134").Value = 3148.2222
$ws.Range("I134").Value = 3148.2222
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9444.6666
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -4374.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 30000
$ws.Range("J57").Value = 30000
$ws.Range("L57").Value = 30000
$ws.Range("N57").Value = -31640

$ws.Range("H132").Value = 3762.125
$ws.Range("I132").Value = 2826.5
$ws.Range("J132").Value = 5820.5
$ws.Range("K132").Value = 8479.5
$ws.Range("L132").Value = 17461.5
$ws.Range("M132").Value = -5949.5
$ws.Range("N132").Value = -22521.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 293.06668
$ws.Range("I55").Value = 253.23077
$ws.Range("J55").Value = 552
$ws.Range("K55").Value = 253.23077
$ws.Range("L55").Value = 552
$ws.Range("M55").Value = -80.23077000000001
$ws.Range("N55").Value = -898

$ws.Range("H68").Value = 2306.6365
$ws.Range("I68").Value = 1961.6666
$ws.Range("J68").Value = 2720.6
$ws.Range("K68").Value = 1961.6666
$ws.Range("L68").Value = 2720.6
$ws.Range("M68").Value = -1212.6666
$ws.Range("N68").Value = -4218.6

$ws.Range("H71").Value = 2306.6365
$ws.Range("I71").Value = 1961.6666
$ws.Range("J71").Value = 2720.6
$ws.Range("K71").Value = 9808.333000000001
$ws.Range("L71").Value = 13603
$ws.Range("M71").Value = -6064.333000000001
$ws.Range("N71").Value = -21091

$ws.Range("H81").Value = 31147.666
$ws.Range("J81").Value = 31147.666
$ws.Range("L81").Value = 31147.666
$ws.Range("N81").Value = -33143.666

$ws.Range("H84").Value = 31147.666
$ws.Range("J84").Value = 31147.666
$ws.Range("L84").Value = 93442.99800000001
$ws.Range("N84").Value = -103426.998

$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52246

$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -161232

$ws.Range("H100").Value = 2731.3447
$ws.Range("I100").Value = 1855.3636
$ws.Range("J100").Value = 3266.6667
$ws.Range("K100").Value = 1855.3636
$ws.Range("L100").Value = 3266.6667
$ws.Range("M100").Value = -1314.3636
$ws.Range("N100").Value = -4348.6667

$ws.Range("H132").Value = 3438.2703
$ws.Range("I132").Value = 2200.0356
$ws.Range("J132").Value = 7290.5557
$ws.Range("K132").Value = 6600.1068
$ws.Range("L132").Value = 21871.6671
$ws.Range("M132").Value = -4070.1068
$ws.Range("N132").Value = -26931.6671

$ws.Range("H136").Value = 4754.7144
$ws.Range("I136").Value = 2049.6875
$ws.Range("J136").Value = 13410.8
$ws.Range("K136").Value = 6149.0625
$ws.Range("L136").Value = 40232.39999999999
$ws.Range("M136").Value = -3599.0625
$ws.Range("N136").Value = -45332.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3644.1035
$ws.Range("I81").Value = 2205.6428
$ws.Range("J81").Value = 4986.6665
$ws.Range("K81").Value = 4411.2856
$ws.Range("L81").Value = 9973.333000000001
$ws.Range("M81").Value = -3350.2856
$ws.Range("N81").Value = -12095.333

$ws.Range("H84").Value = 3644.1035
$ws.Range("I84").Value = 2205.6428
$ws.Range("J84").Value = 4986.6665
$ws.Range("K84").Value = 22056.428
$ws.Range("L84").Value = 49866.665
$ws.Range("M84").Value = -16752.428
$ws.Range("N84").Value = -60474.665

$ws.Range("H132").Value = 2672.4314
$ws.Range("I132").Value = 2473.4614
$ws.Range("J132").Value = 3319.0833
$ws.Range("K132").Value = 7420.3842
$ws.Range("L132").Value = 9957.249899999999
$ws.Range("M132").Value = -4890.3842
$ws.Range("N132").Value = -15017.2499
